# Add a new "2022-Q3" sheet in front of the existing "2022-Q2" sheet,
# populate it with the new quarter's fund-holdings table, and update the
# "总计" (summary) sheet so it lists Q3 first (Q2 / Q1 rows shift down,
# keeping their original counts/values).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Insert the new "2022-Q3" worksheet right before "2022-Q2".
# ---------------------------------------------------------------------
$q2 = $wb.Worksheets.Item("2022-Q2")
$q3 = $wb.Worksheets.Add($q2)
$q3.Name = "2022-Q3"

# Header row
$q3.Range("B1").Value = "基金代码"
$q3.Range("C1").Value = "基金名称"
$q3.Range("D1").Value = "基金规模"
$q3.Range("E1").Value = "股票总仓位"
$q3.Range("F1").Value = "仓位占比"
$q3.Range("G1").Value = "持有市值(亿元)"
$q3.Range("H1").Value = "仓位排名"

# Match the look of the other quarter sheets: bold header, thin border,
# centered / top aligned text, same treatment for the index column (A).
$q3.Range("A1:H7").Borders.LineStyle = 1
$q3.Range("B1:H1").Font.Bold = $true
$q3.Range("A1:H1").HorizontalAlignment = -4108
$q3.Range("A1:H1").VerticalAlignment = -4160
$q3.Range("A2:A7").HorizontalAlignment = -4108
$q3.Range("A2:A7").VerticalAlignment = -4160
$q3.Range("A2:A7").Font.Bold = $true

# Data rows (fund code / name / size / position stay as text like the
# other quarter sheets; rank (H) stays numeric).
$q3Rows = @(
    @(0, "009986", "天弘创新领航混合A",   "2.16", "93.84", "4.07", "0.0879", 10),
    @(1, "012259", "天弘鑫悦成长混合C",   "1.19", "93.52", "4.50", "0.0536", 9),
    @(2, "015769", "天弘低碳经济混合A",   "1.19", "79.16", "3.44", "0.0409", 6),
    @(3, "015770", "天弘低碳经济混合C",   "1.03", "79.16", "3.44", "0.0354", 6),
    @(4, "009987", "天弘创新领航混合C",   "0.53", "93.84", "4.07", "0.0216", 10),
    @(5, "012258", "天弘鑫悦成长混合A",   "0.28", "93.52", "4.50", "0.0126", 9)
)

$r = 2
foreach ($row in $q3Rows) {
    $q3.Range("A$r").Value = $row[0]

    $q3.Range("B$r").NumberFormat = "@"
    $q3.Range("B$r").Value = $row[1]
    $q3.Range("B$r").NumberFormat = "General"

    $q3.Range("C$r").Value = $row[2]

    $q3.Range("D$r").NumberFormat = "@"
    $q3.Range("D$r").Value = $row[3]
    $q3.Range("D$r").NumberFormat = "General"

    $q3.Range("E$r").NumberFormat = "@"
    $q3.Range("E$r").Value = $row[4]
    $q3.Range("E$r").NumberFormat = "General"

    $q3.Range("F$r").NumberFormat = "@"
    $q3.Range("F$r").Value = $row[5]
    $q3.Range("F$r").NumberFormat = "General"

    $q3.Range("G$r").NumberFormat = "@"
    $q3.Range("G$r").Value = $row[6]
    $q3.Range("G$r").NumberFormat = "General"

    $q3.Range("H$r").Value = $row[7]

    $r = $r + 1
}

# ---------------------------------------------------------------------
# 2. Update the "总计" summary sheet: Q3 becomes row 2, Q2 / Q1 shift
#    down a row (their counts/market values are unchanged).
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q3"
$total.Range("C2").Value = 6
$total.Range("D2").Value = 0.25

$total.Range("A3").Value = 1
$total.Range("B3").Value = "2022-Q2"
$total.Range("C3").Value = 6
$total.Range("D3").Value = 0.66

$total.Range("A4").Value = 2
$total.Range("B4").Value = "2022-Q1"
$total.Range("C4").Value = 3
$total.Range("D4").Value = 0.45

# "2022-Q1" was the active tab before this edit (it still is the sheet the
# file was last focused on); keep that selection rather than leaving the
# newly-added "2022-Q3" sheet active.
$wb.Worksheets.Item("2022-Q1").Activate()
